$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "81, 100, 122, 150, 240, 303, 401"
$ws.Range("E3").Value = "12, 6.8, 5.6, 3.1, 1.5, 1.6, 1.3"

$ws.Range("E4").Select()
